# System tests Avanzamento e introduzione nuovi SyT
#
# The "System test" table (rows 21-25) gains a brand-new "Sezione Graph"
# entry as the new row 23; the former rows 23-25 (Smells/Metrics/Members)
# shift down by one and their SyT ids are renumbered, and a new row 26
# ("Sezione Members", SYSTE-18) is appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Materialise the new row 26 by copying row 25's formatting (style,
#        borders, height) down one row, then overwrite its content below. ---
$ws.Range("A25:G25").Copy($ws.Range("A26:G26"))
$ws.Rows.Item(26).RowHeight = 32

# --- 2. Renumber the SyT ids in place, bottom-up, so the workbook's
#        shared-string table grows in the same order the author typed
#        the new identifiers (SYSTE-17, then SYSTE-18, then the new
#        "Sezione Graph" section name). ---
$ws.Range("A25").Value = "SYSTE-17"
$ws.Range("B25").Value = "Sezione Metrics"
$ws.Range("A26").Value = "SYSTE-18"

$ws.Range("A24").Value = "SYSTE-16"
$ws.Range("B24").Value = "Sezione Smells"

$ws.Range("A23").Value = "SYSTE-15"
$ws.Range("B23").Value = "Sezione Graph"

$ws.Range("A22").Value = "SYSTE-14"
$ws.Range("A21").Value = "SYSTE-13"

# --- 3. Fill in the rest of row 26 (same content pattern as the other
#        "Sezione X" rows: Accessibile / Visibile dopo richiesta / ...). ---
$ws.Range("B26").Value = "Sezione Members"
$ws.Range("C26").Value = "Accessibile"
$ws.Range("D26").Value = "Visibile dopo richiesta"
$ws.Range("F26").Value = "Visibile dopo richiesta"
$ws.Range("G26").Value = "Repository: https://github.com/rubygems/bundler" + [char]10 + "Data fine: 01/01/2019"

$excel.CutCopyMode = $false

# --- 4. View tweaks: drop the frozen top-left cell, zoom to 75%, and move
#        the active selection onto the newly edited cell. ---
$excel.ActiveWindow.Zoom = 75
$ws.Range("D23").Select()
